$wb = $excel.ActiveWorkbook

# "Sheet1" worksheet holds the driving value in I23 that feeds the
# CONCATENATE formulas in columns A, B and C (rows 23-42) via $I$23.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("I23").Value = 7

# "order" worksheet stores the same generated strings as static values
# (not formulas) in columns R, S, T for rows 2-21. Update them to match
# the new suffix ("...6" -> "...7").
$order = $wb.Worksheets.Item("order")

$names = @(
    "DonnellJernigan",
    "MalikOtoole",
    "AlanCaudill",
    "AdanApplegate",
    "AiyanaWhitworth",
    "MercedezBrien",
    "DuaneHager",
    "LorenBell",
    "GeraldHiller",
    "DeionBranch",
    "DakotaHalstead",
    "ElliottFurman",
    "MiltonCamp",
    "DawnChester",
    "ZacheryPetrie",
    "EstebanAngel",
    "JimmyBlankenship",
    "AllysaGrice",
    "AugustineYoo",
    "BrandiSouthard"
)

for ($idx = 0; $idx -lt $names.Count; $idx++) {
    $row = 2 + $idx
    $name = $names[$idx]
    $order.Range("R$row").Value = "$name" + "7"
    $order.Range("S$row").Value = "$name" + "7"
    $order.Range("T$row").Value = "$name" + "7@gmail.com"
}
